$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc26fadb063d9f052b50a22571eba7e399e73cdf/e2e/a.md"

# ---------------------------------------------------------------
# Overview sheet: refresh the per-locale status text and widen the
# two status columns (zh-cn / de-de) so the longer text fits.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns("E").ColumnWidth = 29.1
$wsOverview.Columns("F").ColumnWidth = 29.1

# ---------------------------------------------------------------
# zh-cn sheet: status text, widen Status + Latest Handback File
# columns, and record that file a.md has now been handed back
# (target file + handback file + handback datetime + hyperlink).
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText
$wsZh.Columns("C").ColumnWidth = 29.1
$wsZh.Columns("J").ColumnWidth = 39.17

$wsZh.Range("I2").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aUrl, "", "", "a.md")
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-25 14:37:21"

$wsZh.Range("I3").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aUrl, "", "", "a.md")
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-25 14:37:21"

# ---------------------------------------------------------------
# de-de sheet: same treatment, with its own handback datetime.
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText
$wsDe.Columns("C").ColumnWidth = 29.1
$wsDe.Columns("J").ColumnWidth = 39.17

$wsDe.Range("I2").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aUrl, "", "", "a.md")
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-25 14:37:28"

$wsDe.Range("I3").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aUrl, "", "", "a.md")
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-25 14:37:28"
